$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- 1. Rename project in row 39 from "Nottingham and Nottinghamshire" to "Nottingham City" ---
$ws.Range("B39").Value = "Nottingham City"

# --- 2. Add a new row to Table1 for the "Herts and West Essex" project (07H) ---
$tbl = $ws.ListObjects.Item("Table1")
$newRow = $tbl.ListRows.Add()
$r = $newRow.Range

$r.Cells.Item(1, 1).Value = "07H"
$r.Cells.Item(1, 2).Value = "Herts and West Essex"
$r.Cells.Item(1, 2).Font.Name = "Calibri"
$r.Cells.Item(1, 2).VerticalAlignment = -4108
$r.Cells.Item(1, 3).Value = "Phase 3"
$r.Cells.Item(1, 4).Value = 44866
$r.Cells.Item(1, 4).NumberFormat = "mmm-yy"
$r.Cells.Item(1, 5).Value = 800
$r.Cells.Item(1, 6).Value = "16/03/2023 pop supplied by poppy"
$r.Cells.Item(1, 7).Value = "E56000023"
$r.Cells.Item(1, 8).Value = "East of England - South"

# --- 3. Widen column F so the long "status" notes are fully visible ---
$ws.Columns.Item(6).ColumnWidth = 103.65

# --- 4. Update the active selection to reflect where the author left off ---
$ws.Range("F26").Select()
